$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 997.35956
$ws.Range("I15").Value = 997.35956
$ws.Range("K15").Value = 2992.07868
$ws.Range("M15").Value = -2823.07868
$ws.Range("H28").Value = 598.63336
$ws.Range("I28").Value = 626.0476
$ws.Range("J28").Value = 534.6667
$ws.Range("K28").Value = 626.0476
$ws.Range("L28").Value = 534.6667
$ws.Range("M28").Value = -141.0476
$ws.Range("N28").Value = -1504.6667
$ws.Range("H41").Value = 575
$ws.Range("I41").Value = 325
$ws.Range("J41").Value = 675
$ws.Range("K41").Value = 325
$ws.Range("L41").Value = 675
$ws.Range("M41").Value = 115
$ws.Range("N41").Value = -1555
$ws.Range("H53").Value = 436.10715
$ws.Range("I53").Value = 322.3846
$ws.Range("J53").Value = 534.6667
$ws.Range("K53").Value = 322.3846
$ws.Range("L53").Value = 534.6667
$ws.Range("M53").Value = 314.6154
$ws.Range("N53").Value = -1808.6667
$ws.Range("H76").Value = 3021.6956
$ws.Range("I76").Value = 2977.2273
$ws.Range("J76").Value = 4000
$ws.Range("K76").Value = 2977.2273
$ws.Range("L76").Value = 4000
$ws.Range("M76").Value = -2662.2273
$ws.Range("N76").Value = -4630
$ws.Range("H79").Value = 3021.6956
$ws.Range("I79").Value = 2977.2273
$ws.Range("J79").Value = 4000
$ws.Range("K79").Value = 2977.2273
$ws.Range("L79").Value = 4000
$ws.Range("M79").Value = -1885.2273
$ws.Range("N79").Value = -6184
$ws.Range("H92").Value = 1664.3158
$ws.Range("I92").Value = 1647.7142
$ws.Range("J92").Value = 1710.8
$ws.Range("K92").Value = 1647.7142
$ws.Range("L92").Value = 1710.8
$ws.Range("M92").Value = -399.7141999999999
$ws.Range("N92").Value = -4206.8
$ws.Range("H98").Value = 2520.75
$ws.Range("I98").Value = 966.3684
$ws.Range("J98").Value = 4792.5386
$ws.Range("K98").Value = 966.3684
$ws.Range("L98").Value = 4792.5386
$ws.Range("M98").Value = 531.6316
$ws.Range("N98").Value = -7788.5386
$ws.Range("H112").Value = 1324.6792
$ws.Range("I112").Value = 831.2857
$ws.Range("K112").Value = 2493.8571
$ws.Range("M112").Value = -1385.8571
$ws.Range("H122").Value = 2520.75
$ws.Range("I122").Value = 966.3684
$ws.Range("J122").Value = 4792.5386
$ws.Range("K122").Value = 2899.1052
$ws.Range("L122").Value = 14377.6158
$ws.Range("M122").Value = -449.1052
$ws.Range("N122").Value = -19277.6158
$ws.Range("H132").Value = 184903.77
$ws.Range("I132").Value = 2957.4893
$ws.Range("J132").Value = 1253838.1
$ws.Range("K132").Value = 8872.4679
$ws.Range("L132").Value = 3761514.3
$ws.Range("M132").Value = -6342.4679
$ws.Range("N132").Value = -3766574.3
$ws.Range("H135").Value = 250.76666
$ws.Range("I135").Value = 205.95833
$ws.Range("J135").Value = 430
$ws.Range("K135").Value = 1853.62497
$ws.Range("L135").Value = 3870
$ws.Range("M135").Value = 681.3750300000002
$ws.Range("N135").Value = -8940

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5232.528
$ws.Range("I32").Value = 4134.4194
$ws.Range("K32").Value = 4134.4194
$ws.Range("M32").Value = -3847.4194
$ws.Range("H94").Value = 34975
$ws.Range("J94").Value = 34975
$ws.Range("L94").Value = 34975
$ws.Range("N94").Value = -36777
$ws.Range("H109").Value = 31000
$ws.Range("J109").Value = 31000
$ws.Range("L109").Value = 31000
$ws.Range("N109").Value = -33774

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13163.6
$ws.Range("I20").Value = 5999.1665
$ws.Range("J20").Value = 17939.889
$ws.Range("K20").Value = 5999.1665
$ws.Range("L20").Value = 17939.889
$ws.Range("M20").Value = -5752.1665
$ws.Range("N20").Value = -18433.889
$ws.Range("H134").Value = 1761.0244
$ws.Range("I134").Value = 1218.0625
$ws.Range("J134").Value = 3691.5557
$ws.Range("K134").Value = 3654.1875
$ws.Range("L134").Value = 11074.6671
$ws.Range("M134").Value = -1119.1875
$ws.Range("N134").Value = -16144.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1400.909
$ws.Range("I105").Value = 1275.3334
$ws.Range("J105").Value = 1966
$ws.Range("K105").Value = 1275.3334
$ws.Range("L105").Value = 1966
$ws.Range("M105").Value = 471.6666
$ws.Range("N105").Value = -5460

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 209.06667
$ws.Range("J23").Value = 242.16667
$ws.Range("L23").Value = 726.50001
$ws.Range("N23").Value = -1196.50001
$ws.Range("H34").Value = 13108.353
$ws.Range("I34").Value = 36731.332
$ws.Range("J34").Value = 8046.2856
$ws.Range("K34").Value = 110193.996
$ws.Range("L34").Value = 24138.8568
$ws.Range("M34").Value = -110109.996
$ws.Range("N34").Value = -24306.8568
$ws.Range("H39").Value = 8258.786
$ws.Range("J39").Value = 9602.091
$ws.Range("L39").Value = 28806.273
$ws.Range("N39").Value = -29394.273
$ws.Range("L75").Value = 5850
$ws.Range("M75").Value = -202
$ws.Range("N75").Value = -7846
$ws.Range("H78").Value = 1433.3334
$ws.Range("I78").Value = 400
$ws.Range("J78").Value = 1950
$ws.Range("K78").Value = 3600
$ws.Range("L78").Value = 17550
$ws.Range("M78").Value = 1392
$ws.Range("N78").Value = -27534
$ws.Range("H117").Value = 4069.375
$ws.Range("J117").Value = 30000
$ws.Range("L117").Value = 90000
$ws.Range("N117").Value = -96884

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 23407
$ws.Range("J93").Value = 23407
$ws.Range("L93").Value = 23407
$ws.Range("N93").Value = -27151
$ws.Range("H102").Value = 3144.5
$ws.Range("I102").Value = 2166.75
$ws.Range("J102").Value = 5100
$ws.Range("K102").Value = 2166.75
$ws.Range("L102").Value = 5100
$ws.Range("M102").Value = -544.75
$ws.Range("N102").Value = -8344
$ws.Range("H126").Value = 2021.37
$ws.Range("I126").Value = 2038.9166
$ws.Range("J126").Value = 1600.25
$ws.Range("K126").Value = 6116.7498
$ws.Range("L126").Value = 4800.75
$ws.Range("M126").Value = -3646.7498
$ws.Range("N126").Value = -9740.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1911.4706
$ws.Range("I100").Value = 1722.3846
$ws.Range("J100").Value = 2526
$ws.Range("K100").Value = 1722.3846
$ws.Range("L100").Value = 2526
$ws.Range("M100").Value = -1181.3846
$ws.Range("N100").Value = -3608
$ws.Range("H132").Value = 3765.524
$ws.Range("I132").Value = 1709.7307
$ws.Range("J132").Value = 7106.1875
$ws.Range("K132").Value = 5129.1921
$ws.Range("L132").Value = 21318.5625
$ws.Range("M132").Value = -2599.1921
$ws.Range("N132").Value = -26378.5625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 465.66666
$ws.Range("I107").Value = 398.5
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 1195.5
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 724.5
$ws.Range("N107").Value = -5640
$ws.Range("H122").Value = 3523.0527
$ws.Range("I122").Value = 1996.125
$ws.Range("J122").Value = 11666.667
$ws.Range("K122").Value = 5988.375
$ws.Range("L122").Value = 35000.001
$ws.Range("M122").Value = -3538.375
$ws.Range("N122").Value = -39900.001
$ws.Range("H139").Value = 37513.684
$ws.Range("J139").Value = 37339.445
$ws.Range("L139").Value = 37339.445
$ws.Range("N139").Value = -47619.445
